$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.945.93'
$ws.Range('E2').Value = '  -3.24%  '
$ws.Range('D3').Value = '3.227.67'
$ws.Range('E3').Value = '  -3.88%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '''537.40'
$ws.Range('E5').Value = '  -5.48%  '
$ws.Range('D6').Value = '''135.98'
$ws.Range('E6').Value = '  -9.48%  '
$ws.Range('E7').Value = '  -0.15%  '
$ws.Range('D8').Value = '3.227.79'
$ws.Range('E8').Value = '  -4.03%  '
$ws.Range('D9').Value = '''0.458'
$ws.Range('E9').Value = '  -4.61%  '
$ws.Range('D10').Value = '''7.60'
$ws.Range('E10').Value = '  -4.15%  '
$ws.Range('E11').Value = '  -4.92%  '
$ws.Range('E12').Value = '  -4.49%  '
$ws.Range('D13').Value = '3.784.32'
$ws.Range('E13').Value = '  -4.32%  '
$ws.Range('E14').Value = '  -0.99%  '
$ws.Range('D15').Value = '''26.14'
$ws.Range('E15').Value = '  -7.00%  '
$ws.Range('D16').Value = '3.227.22'
$ws.Range('E16').Value = '  -4.15%  '
$ws.Range('D17').Value = '''0.0000159'
$ws.Range('E17').Value = '  -5.95%  '
$ws.Range('D18').Value = '59.062.87'
$ws.Range('E18').Value = '  -3.29%  '
$ws.Range('D19').Value = '''5.91'
$ws.Range('E19').Value = '  -6.52%  '
$ws.Range('D20').Value = '''13.24'
$ws.Range('E20').Value = '  -6.48%  '
$ws.Range('D21').Value = '''8.29'
$ws.Range('E21').Value = '  -5.87%  '
$ws.Range('D22').Value = '''361.30'
$ws.Range('E22').Value = '  -2.95%  '
$ws.Range('E23').Value = '  +0.01%  '
$ws.Range('D24').Value = '''70.60'
$ws.Range('E24').Value = '  -6.16%  '
$ws.Range('D25').Value = '''0.519'
$ws.Range('E25').Value = '  -7.71%  '
$ws.Range('D26').Value = '3.361.65'
$ws.Range('E26').Value = '  -4.96%  '
$ws.Range('B27').Value = 'Kaspa'
$ws.Range('C27').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D27').Value = '''0.171'
$ws.Range('E27').Value = '  -1.90%  '
$ws.Range('B28').Value = 'PEPE'
$ws.Range('C28').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D28').Value = '0.0₃0977'
$ws.Range('E28').Value = '  -9.14%  '
$ws.Range('E29').Value = '  -0.06%  '
$ws.Range('E30').Value = '  -4.44%  '
$ws.Range('E31').Value = '  -0.17%  '
$ws.Range('E32').Value = '  -7.33%  '
$ws.Range('D33').Value = '''7.11'
$ws.Range('E33').Value = '  -7.28%  '
$ws.Range('D34').Value = '''22.00'
$ws.Range('E34').Value = '  -3.75%  '
$ws.Range('E35').Value = '  -1.87%  '
$ws.Range('B36').Value = 'Monero'
$ws.Range('C36').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D36').Value = '''163.21'
$ws.Range('E36').Value = '  -3.75%  '
$ws.Range('B37').Value = 'NEARProtocol'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D37').Value = '''4.95'
$ws.Range('E37').Value = '  -7.41%  '
$ws.Range('D38').Value = '''6.37'
$ws.Range('E38').Value = '  -5.45%  '
$ws.Range('D39').Value = '''1.44'
$ws.Range('E39').Value = '  -6.27%  '
$ws.Range('D40').Value = '''25.84'
$ws.Range('E40').Value = '  -14.08%  '
$ws.Range('D41').Value = '''0.0707'
$ws.Range('E41').Value = '  -6.23%  '
$ws.Range('D42').Value = '3.260.96'
$ws.Range('E42').Value = '  -4.32%  '
$ws.Range('D43').Value = '''41.02'
$ws.Range('E43').Value = '  -3.18%  '
$ws.Range('D44').Value = '''0.718'
$ws.Range('E44').Value = '  -5.86%  '
$ws.Range('E45').Value = '  -6.64%  '
$ws.Range('E46').Value = '  -3.88%  '
$ws.Range('D47').Value = '''1.50'
$ws.Range('E47').Value = '  -6.16%  '
$ws.Range('E48').Value = '  -0.25%  '
$ws.Range('D49').Value = '2.291.96'
$ws.Range('E49').Value = '  -8.63%  '
$ws.Range('D50').Value = '''6.26'
$ws.Range('E50').Value = '  -6.36%  '
$ws.Range('E51').Value = '  -8.77%  '
